$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Select rows 10 through 13 (entire rows) and delete them, shifting rows 14-18 up.
$rng = $ws.Range("A10:XFD13")
$rng.EntireRow.Select()
$rng.EntireRow.Delete()
